$wb = $excel.ActiveWorkbook

# Update the "想去人数" (interest count) figures on both the "展览" and
# "全部类型" sheets, which carry duplicate rows for the same events.
foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 213
    $ws.Range("F4").Value = 3605
    $ws.Range("F5").Value = 374
}
